$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.007.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "'1.862.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("D4").Value = "'0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'305.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").Value = "'0.9981"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5063"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.07141"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'20.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.8833"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07564"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.843.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "'5.304"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "'89.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "'0.9989"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'0.000008403"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "'0.9978"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'27.043.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'2.112.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'10.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").Value = "'6.460"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "'147.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("D27").Value = "'17.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").Value = "'2.102"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").Value = "'112.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'4.677"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").Value = "'4.704"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "'0.05137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "'3.034"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").Value = "'1.152"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.85%  "
$ws.Range("D36").Value = "'0.7287"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("D37").Value = "'0.02038"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'3.033"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "'2.466"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.73%  "
$ws.Range("D40").Value = "'1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "'0.5290"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").Value = "'6.548"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").Value = "'116.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "'8.274"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").Value = "'0.9974"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'0.4617"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").Value = "'10.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "'1.564"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").Value = "'36.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "'63.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.81%  "
